$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value. All target cells are plain
# text (coin names/links/prices/volume %), so we force the Text number
# format before the write (otherwise Excel auto-coerces number-looking
# strings like "1.011" or "20.551.38" into actual numbers), then restore
# the default "Normal" cell style so no stray formatting is left behind.
$updates = @(
    @{Cell="D2"; Value="20.551.38"}
    @{Cell="E2"; Value="  +0.05%  "}
    @{Cell="D3"; Value="1.477.89"}
    @{Cell="E3"; Value="  +0.67%  "}
    @{Cell="D4"; Value="1.011"}
    @{Cell="E4"; Value="  +0.06%  "}
    @{Cell="D5"; Value="0.9754"}
    @{Cell="E5"; Value="  +2.45%  "}
    @{Cell="D6"; Value="279.08"}
    @{Cell="E6"; Value="  -0.87%  "}
    @{Cell="E7"; Value="  -1.20%  "}
    @{Cell="D8"; Value="0.3073"}
    @{Cell="E8"; Value="  -3.59%  "}
    @{Cell="D9"; Value="40.01"}
    @{Cell="E9"; Value="  -4.49%  "}
    @{Cell="D10"; Value="1.059"}
    @{Cell="E10"; Value="  +0.02%  "}
    @{Cell="D11"; Value="0.06665"}
    @{Cell="E11"; Value="  -0.26%  "}
    @{Cell="D12"; Value="1.005"}
    @{Cell="E12"; Value="  -0.03%  "}
    @{Cell="D13"; Value="5.497"}
    @{Cell="E13"; Value="  -2.29%  "}
    @{Cell="D14"; Value="18.03"}
    @{Cell="E14"; Value="  -1.02%  "}
    @{Cell="D15"; Value="6.201"}
    @{Cell="E15"; Value="  -1.05%  "}
    @{Cell="D16"; Value="0.9771"}
    @{Cell="E16"; Value="  +2.88%  "}
    @{Cell="D17"; Value="0.00001027"}
    @{Cell="E17"; Value="  -0.68%  "}
    @{Cell="D18"; Value="1.479.41"}
    @{Cell="E18"; Value="  +0.18%  "}
    @{Cell="D19"; Value="0.05932"}
    @{Cell="E19"; Value="  +4.65%  "}
    @{Cell="D20"; Value="69.36"}
    @{Cell="E20"; Value="  -4.02%  "}
    @{Cell="D21"; Value="5.479"}
    @{Cell="E21"; Value="  -3.68%  "}
    @{Cell="D22"; Value="14.49"}
    @{Cell="E22"; Value="  -1.34%  "}
    @{Cell="D23"; Value="11.04"}
    @{Cell="E23"; Value="  -1.54%  "}
    @{Cell="D24"; Value="2.255"}
    @{Cell="E24"; Value="  -1.16%  "}
    @{Cell="D25"; Value="20.623.95"}
    @{Cell="E25"; Value="  -0.20%  "}
    @{Cell="D26"; Value="142.11"}
    @{Cell="E26"; Value="  +3.24%  "}
    @{Cell="D27"; Value="2.127"}
    @{Cell="E27"; Value="  -7.52%  "}
    @{Cell="E28"; Value="  -1.85%  "}
    @{Cell="D29"; Value="1.640.15"}
    @{Cell="E29"; Value="  +0.07%  "}
    @{Cell="D30"; Value="113.86"}
    @{Cell="E30"; Value="  -0.01%  "}
    @{Cell="D31"; Value="3.932"}
    @{Cell="E31"; Value="  -0.59%  "}
    @{Cell="D32"; Value="5.007"}
    @{Cell="E32"; Value="  -5.85%  "}
    @{Cell="D33"; Value="0.8154"}
    @{Cell="E33"; Value="  -2.82%  "}
    @{Cell="D34"; Value="0.07996"}
    @{Cell="E34"; Value="  +1.77%  "}
    @{Cell="D35"; Value="1.527"}
    @{Cell="E35"; Value="  -6.47%  "}
    @{Cell="D36"; Value="1.217"}
    @{Cell="E36"; Value="  +8.64%  "}
    @{Cell="D37"; Value="0.05783"}
    @{Cell="E37"; Value="  -4.51%  "}
    @{Cell="D38"; Value="4.709"}
    @{Cell="E38"; Value="  -4.28%  "}
    @{Cell="B39"; Value="Frax"}
    @{Cell="C39"; Value="https://coinranking.com/coin/KfWtaeV1W+frax-frax"}
    @{Cell="D39"; Value="0.9762"}
    @{Cell="E39"; Value="  +1.27%  "}
    @{Cell="B40"; Value="FraxShare"}
    @{Cell="C40"; Value="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"}
    @{Cell="D40"; Value="7.696"}
    @{Cell="E40"; Value="  +3.67%  "}
    @{Cell="D41"; Value="0.02042"}
    @{Cell="E41"; Value="  -1.29%  "}
    @{Cell="D42"; Value="10.43"}
    @{Cell="E42"; Value="  -1.89%  "}
    @{Cell="D43"; Value="0.1887"}
    @{Cell="E43"; Value="  -0.14%  "}
    @{Cell="D44"; Value="0.5289"}
    @{Cell="E44"; Value="  -2.47%  "}
    @{Cell="D45"; Value="3.529"}
    @{Cell="E45"; Value="  -1.71%  "}
    @{Cell="D46"; Value="12.23"}
    @{Cell="E46"; Value="  -2.90%  "}
    @{Cell="D47"; Value="118.56"}
    @{Cell="E47"; Value="  -3.16%  "}
    @{Cell="D48"; Value="0.5193"}
    @{Cell="E48"; Value="  -2.69%  "}
    @{Cell="D49"; Value="1.802"}
    @{Cell="E49"; Value="  -1.78%  "}
    @{Cell="E50"; Value="  +0.65%  "}
    @{Cell="D51"; Value="0.9962"}
    @{Cell="E51"; Value="  -0.01%  "}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
